# The commit deletes the content slide titled "Gabarito - Exercicios
# Selecionados" (the 8th slide, "Slide 8" / ppt/slides/slide8.xml). Removing
# it through the Slides collection also drops its notes page and lets the
# host re-number the surrounding p:sldId / relationship bookkeeping the way
# PowerPoint itself does on save.
$p = $ppt.ActivePresentation
$p.Slides.Item(8).Delete()
